# Fill in the "2019" rows that were left blank for Concessionária (col A)
# and Ano (col B) in each of the per-company blocks of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Rows where A/B cells already exist (blank) with the right formatting
# (style s="15" on col A, s="18" on col B) - just fill in the values.
$updates = @(
    @{ Row = 30;  Name = "EFVM" },
    @{ Row = 45;  Name = "FCA" },
    @{ Row = 60;  Name = "EFPO" },
    @{ Row = 73;  Name = "FNSTN" },
    @{ Row = 88;  Name = "FTC" },
    @{ Row = 103; Name = "FTL" },
    @{ Row = 118; Name = "MRS" },
    @{ Row = 133; Name = "RMN" },
    @{ Row = 148; Name = "RMO" },
    @{ Row = 163; Name = "RMP" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Name
    $ws.Cells.Item($u.Row, 2).Value = 2019
}

# Row 178: A178/B178 don't exist yet in the sheet (the row only has D178),
# so pick up the formatting of the row above (A177/B177, same style indices
# 15/18 used throughout this column) before writing the new values.
$ws.Range("A177:B177").Copy()
$ws.Range("A178:B178").PasteSpecial(-4122)
$ws.Cells.Item(178, 1).Value = "RMS"
$ws.Cells.Item(178, 2).Value = 2019

$excel.CutCopyMode = 0

# Leave the selection where the user ended up after filling in the last
# block (RMS / 2019 in A178:B178).
$ws.Activate() | Out-Null
$ws.Range("A177:B178").Select() | Out-Null
